$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.860.01"
$ws.Range("E2").Value = "  +3.92%  "

$ws.Range("D3").Value = "2.486.08"
$ws.Range("E3").Value = "  +2.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("E7").Value = "  +0.79%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.60%  "

$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("E12").Value = "  +0.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.27%  "

$ws.Range("D15").Value = "2.872.05"
$ws.Range("E15").Value = "  +1.92%  "

$ws.Range("D16").Value = "2.532.61"
$ws.Range("E16").Value = "  +3.94%  "

$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("D18").Value = "46.797.89"
$ws.Range("E18").Value = "  +3.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.18%  "

$ws.Range("D21").Value = "0.0₃0931"
$ws.Range("E21").Value = "  +0.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.78%  "

$ws.Range("E24").Value = "  +2.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("E28").Value = "  +4.32%  "

$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.132"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("E33").Value = "  -2.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.10%  "

$ws.Range("E35").Value = "  +1.09%  "

$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("E37").Value = "  +0.69%  "

$ws.Range("E38").Value = "  +0.98%  "

$ws.Range("E39").Value = "  +2.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "122.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.01%  "

$ws.Range("E41").Value = "  +1.08%  "

$ws.Range("E42").Value = "  +1.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "

$ws.Range("E44").Value = "  +1.55%  "

$ws.Range("D45").Value = "1.953.91"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.61%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.28%  "

$ws.Range("E48").Value = "  -1.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.74%  "

$ws.Range("E50").Value = "  +14.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.43%  "
